{"js": "// Replace the date line and the 25 \"a\u00f7b=c, d\" answers in the worksheet\n// table, in document order. Each row of the pair list is\n// [oldText, newText]; old texts are looked up against the body's\n// paragraphs (this also walks into the table cells) so we never rely on\n// absolute table/row/column indices.\nconst pairs = [\n  [\"2025-05-17 Saturday\", \"2025-05-18 Sunday\"],\n  [\"50\u00f79=5, 5\", \"22\u00f77=3, 1\"],\n  [\"62\u00f75=12, 2\", \"42\u00f72=21, 0\"],\n  [\"16\u00f76=2, 4\", \"22\u00f74=5, 2\"],\n  [\"24\u00f76=4, 0\", \"86\u00f75=17, 1\"],\n  [\"71\u00f77=10, 1\", \"35\u00f75=7, 0\"],\n  [\"24\u00f78=3, 0\", \"15\u00f74=3, 3\"],\n  [\"68\u00f78=8, 4\", \"35\u00f76=5, 5\"],\n  [\"13\u00f75=2, 3\", \"29\u00f78=3, 5\"],\n  [\"95\u00f79=10, 5\", \"33\u00f73=11, 0\"],\n  [\"72\u00f75=14, 2\", \"83\u00f74=20, 3\"],\n  [\"61\u00f72=30, 1\", \"60\u00f77=8, 4\"],\n  [\"73\u00f76=12, 1\", \"10\u00f74=2, 2\"],\n  [\"76\u00f77=10, 6\", \"99\u00f78=12, 3\"],\n  [\"10\u00f75=2, 0\", \"11\u00f74=2, 3\"],\n  [\"42\u00f73=14, 0\", \"68\u00f79=7, 5\"],\n  [\"83\u00f73=27, 2\", \"23\u00f78=2, 7\"],\n  [\"75\u00f75=15, 0\", \"13\u00f74=3, 1\"],\n  [\"82\u00f79=9, 1\", \"12\u00f72=6, 0\"],\n  [\"78\u00f79=8, 6\", \"28\u00f75=5, 3\"],\n  [\"91\u00f78=11, 3\", \"68\u00f79=7, 5\"],\n  [\"38\u00f78=4, 6\", \"47\u00f77=6, 5\"],\n  [\"33\u00f73=11, 0\", \"34\u00f76=5, 4\"],\n  [\"37\u00f77=5, 2\", \"43\u00f79=4, 7\"],\n  [\"84\u00f78=10, 4\", \"40\u00f78=5, 0\"],\n  [\"58\u00f73=19, 1\", \"59\u00f72=29, 1\"],\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet pairIndex = 0;\nfor (let i = 0; i < paragraphs.items.length && pairIndex < pairs.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text;\n  if (text === \"\") {\n    continue;\n  }\n  const [oldText, newText] = pairs[pairIndex];\n  if (text !== oldText) {\n    throw new Error(\n      `Unexpected paragraph text at index ${i}: expected \"${oldText}\", found \"${text}\"`\n    );\n  }\n  para.insertText(newText, \"Replace\");\n  pairIndex++;\n}\nawait context.sync();\n\nif (pairIndex !== pairs.length) {\n  throw new Error(`Only replaced ${pairIndex} of ${pairs.length} expected paragraphs`);\n}\n", "ps1": "# Replace the date line and the 25 \"a\u00f7b=c, d\" answers in the worksheet\n# table, in document order. Each entry in $pairs is the [old, new] text\n# for one (non-empty) paragraph; walking $d.Paragraphs() in order and\n# skipping the blank filler rows lines these up with the table cells\n# without depending on absolute table/row/column indices. (A plain\n# document-wide Find/Replace-all is unsafe here: one answer's NEW text\n# is identical to a LATER answer's OLD text, so a global replace would\n# re-match its own freshly written output.)\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @{Old=\"2025-05-17 Saturday\"; New=\"2025-05-18 Sunday\"},\n    @{Old=\"50\u00f79=5, 5\"; New=\"22\u00f77=3, 1\"},\n    @{Old=\"62\u00f75=12, 2\"; New=\"42\u00f72=21, 0\"},\n    @{Old=\"16\u00f76=2, 4\"; New=\"22\u00f74=5, 2\"},\n    @{Old=\"24\u00f76=4, 0\"; New=\"86\u00f75=17, 1\"},\n    @{Old=\"71\u00f77=10, 1\"; New=\"35\u00f75=7, 0\"},\n    @{Old=\"24\u00f78=3, 0\"; New=\"15\u00f74=3, 3\"},\n    @{Old=\"68\u00f78=8, 4\"; New=\"35\u00f76=5, 5\"},\n    @{Old=\"13\u00f75=2, 3\"; New=\"29\u00f78=3, 5\"},\n    @{Old=\"95\u00f79=10, 5\"; New=\"33\u00f73=11, 0\"},\n    @{Old=\"72\u00f75=14, 2\"; New=\"83\u00f74=20, 3\"},\n    @{Old=\"61\u00f72=30, 1\"; New=\"60\u00f77=8, 4\"},\n    @{Old=\"73\u00f76=12, 1\"; New=\"10\u00f74=2, 2\"},\n    @{Old=\"76\u00f77=10, 6\"; New=\"99\u00f78=12, 3\"},\n    @{Old=\"10\u00f75=2, 0\"; New=\"11\u00f74=2, 3\"},\n    @{Old=\"42\u00f73=14, 0\"; New=\"68\u00f79=7, 5\"},\n    @{Old=\"83\u00f73=27, 2\"; New=\"23\u00f78=2, 7\"},\n    @{Old=\"75\u00f75=15, 0\"; New=\"13\u00f74=3, 1\"},\n    @{Old=\"82\u00f79=9, 1\"; New=\"12\u00f72=6, 0\"},\n    @{Old=\"78\u00f79=8, 6\"; New=\"28\u00f75=5, 3\"},\n    @{Old=\"91\u00f78=11, 3\"; New=\"68\u00f79=7, 5\"},\n    @{Old=\"38\u00f78=4, 6\"; New=\"47\u00f77=6, 5\"},\n    @{Old=\"33\u00f73=11, 0\"; New=\"34\u00f76=5, 4\"},\n    @{Old=\"37\u00f77=5, 2\"; New=\"43\u00f79=4, 7\"},\n    @{Old=\"84\u00f78=10, 4\"; New=\"40\u00f78=5, 0\"},\n    @{Old=\"58\u00f73=19, 1\"; New=\"59\u00f72=29, 1\"}\n)\n\n$pairIndex = 0\nforeach ($p in $d.Paragraphs) {\n    if ($pairIndex -ge $pairs.Count) {\n        break\n    }\n    $r = $p.Range\n    $t = $r.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq \"\") {\n        continue\n    }\n    $pair = $pairs[$pairIndex]\n    if ($t -ne $pair.Old) {\n        throw \"Unexpected paragraph text: expected '$($pair.Old)', found '$t'\"\n    }\n    $r.Text = $pair.New\n    $pairIndex++\n}\n\nif ($pairIndex -ne $pairs.Count) {\n    throw \"Only replaced $pairIndex of $($pairs.Count) expected paragraphs\"\n}\n\n$d.Save()\n"}
